# Commit: "added harvester and experiment design"
#
# Changes:
#   - Column B (harvester) value changed from "Retrofitted_0674" to "S.GISH"
#     for every data row (2-18).
#   - Column D (experimentDesign) gets a new value "90minuteInduction" for
#     every data row (2-18), styled like the other data cells (plain,
#     non-bold, black Calibri 11 text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 2).Value = "S.GISH"
    $ws.Cells.Item($row, 4).Value = "90minuteInduction"
}

# New experimentDesign column gets its own (non-bold, explicit black) font.
$ws.Range("D2:D18").Font.Color = 0

# Match the active selection recorded in the edited workbook.
$ws.Range("D2:D18").Select()
